$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "footer" task row (row 18) - the footer work item was dropped.
# This shifts all rows below (header, Homepage realisatie, ... Links realisatie) up by one.
$ws.Rows("18").Delete()

# Mark completed tasks (Klaar? column = F) with an "X" for the mockup rows.
$ws.Range("F8").Value = "X"
$ws.Range("F11").Value = "X"
$ws.Range("F12").Value = "X"
$ws.Range("F13").Value = "X"
$ws.Range("F14").Value = "X"
$ws.Range("F15").Value = "X"
$ws.Range("F16").Value = "X"

# Update the view state to match where the user left off editing.
$ws.Range("E18").Select()
$excel.ActiveWindow.ScrollRow = 10
